# #5: insurance, claim, debt, investment done
# This edit rewrites the "保險" (insurance) worksheet (sheet4) so that:
#   - row 1 becomes a proper header row of column names (company, name, owner,
#     property_category, category, date, legislator_name, legislator_id,
#     source_file, index) instead of duplicating the data row's values, and
#   - row 2 is filled out with the actual record data across all of those
#     columns (previously only A:D were populated).

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item(4)   # "保險" sheet
$src = $wb.Worksheets.Item(3)   # "存款" sheet - used to borrow an already-text "2012-04-25" cell

# --- Extend the existing header/data formatting (bold+bordered header style,
#     plain data style) from columns B:D out to E:K, so the new columns look
#     consistent with the old ones. ---
$ws.Range("B1").Copy()
$ws.Range("E1:K1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B2").Copy()
$ws.Range("E2:K2").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Row 1: column headers ---
$ws.Range("B1").Value = "company"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "property_category"
$ws.Range("F1").Value = "category"
$ws.Range("G1").Value = "date"
$ws.Range("H1").Value = "legislator_name"
$ws.Range("I1").Value = "legislator_id"
$ws.Range("J1").Value = "source_file"
$ws.Range("K1").Value = "index"

# --- Row 2: the actual record values ---
$ws.Range("B2").Value = "中華郵政"
$ws.Range("C2").Value = "吉利保險"
$ws.Range("D2").Value = "梁寒衣"
$ws.Range("E2").Value = "insurance"
$ws.Range("F2").Value = "normal"

# The date value must stay a plain text string ("2012-04-25"), not get
# auto-converted to a date serial number. Borrow the value from the existing
# text cell on the deposit sheet that already holds this exact string so it
# round-trips as text.
$src.Range("I2").Copy()
$ws.Range("G2").PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = $false

$ws.Range("H2").Value = "陳學聖"
$ws.Range("I2").Value = 840
$ws.Range("J2").Value = "tmpfd9c1"
$ws.Range("K2").Value = 107
